$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.287629211722731,0.3402632080236003,0.2352893041228299,0.1776475392211339,1.091381633918587,0.07973214163530429,0.4404514370260948,0.1793734731068,2.400617408154602)
    3 = @(1.142681825063619,0.2985870368923713,0.2286407175794523,0.1732199729132518,1.095427471421289,0.07973214163530429,0.4510568956887298,0.1753624964212861,2.424267857758338)
    4 = @(1.053468655673839,0.2729121377405761,0.2246253972989365,0.170583346297434,1.098777232410193,0.07973214163530429,0.4580034490540967,0.1730093420532199,2.440971406002078)
    5 = @(1.017061812460497,0.2624285278000968,0.2230060757803329,0.1695295362741973,1.100359546011774,0.07973214163530429,0.4609432486398308,0.1720779490771633,2.448325514310483)
    6 = @(1.011013425886119,0.2606864900856465,0.2227382162644602,0.1693557994635881,1.10063539805423,0.07973214163530429,0.4614379759974749,0.1719249544422539,2.449579670792303)
    7 = @(1.052977866751974,0.2727708358038683,0.2246034897529938,0.1705690506382993,1.098797692910402,0.07973214163530429,0.4580426551914396,0.1729966694981684,2.441068372074113)
    8 = @(1.237697270622846,0.3259113889524201,0.2329830487545195,0.1761039334155825,1.092596738313475,0.07973214163530429,0.4440177853716278,0.1779677245160087,2.408318335605969)
    9 = @(1.598144478682229,0.4294182518974594,0.2499421218838052,0.1876068884387436,1.087324285797692,0.07973214163530429,0.4199777512684495,0.1885874164418482,2.361474758672415)
    10 = @(1.861792553738724,0.505014161003885,0.2627187381344669,0.1964537076048956,1.087678212699956,0.07973214163530429,0.4044438692803904,0.196924606132626,2.337745632850073)
    11 = @(1.981462728161716,0.5393025812857104,0.2685991201592515,0.2005643137564093,1.088763494484311,0.07973214163530429,0.39784295483555,0.2008344094956271,2.329291069413074)
    12 = @(2.026738800614737,0.5522717415474858,0.2708355807727116,0.2021332615412774,1.089307848763852,0.07973214163530429,0.3954106240007711,0.2023318455444922,2.326427615556014)
    13 = @(2.016989623658105,0.5494792837974956,0.2703534900808648,0.2017948118586546,1.089184672650958,0.07973214163530429,0.39593147155305,0.2020085947674346,2.327029248373151)
    14 = @(1.985188441542391,0.5403698687111955,0.2687829216774844,0.2006931447238784,1.088805602975896,0.07973214163530429,0.3976414948729623,0.2009572660269612,2.329048705895303)
    15 = @(1.965703964868112,0.5347881002741701,0.267822161105471,0.2000199490965358,1.088590796051108,0.07973214163530429,0.3986977073523885,0.2003154957384083,2.330329761496472)
    16 = @(1.85396629775579,0.5027712521620629,0.262335803737372,0.1961868000946083,1.0876259261629,0.07973214163530429,0.40488465207196,0.1966714515954635,2.33834538300448)
    17 = @(1.785349379317552,0.4831037285050002,0.2589874889391695,0.1938573244432362,1.087271075529245,0.07973214163530429,0.4087996238564138,0.1944659801298201,2.3438631787308)
    18 = @(1.745857971651958,0.4717820651452485,0.2570680568396142,0.1925255821061285,1.087153944196444,0.07973214163530429,0.4110952086980451,0.1932084785110533,2.347257017635343)
    19 = @(1.732482677737721,0.4679471462209221,0.2564192782818964,0.1920760711397591,1.087129208231261,0.07973214163530429,0.4118799657276764,0.1927846031289846,2.34844387434336)
    20 = @(1.792656348374635,0.4851983501445147,0.2593432581890909,0.1941044619716763,1.087299845560025,0.07973214163530429,0.4083783331028084,0.1946996148461295,2.343253003588757)
    21 = @(1.994530335469449,0.5430459403470422,0.2692439737592167,0.2010163961711058,1.088913321237627,0.07973214163530429,0.3971373904914568,0.2012656084578168,2.328446353231499)
    22 = @(2.126229946432659,0.5807642492007972,0.2757710730726046,0.205605715696656,1.090745489015333,0.07973214163530429,0.3901831775703073,0.2056552583962628,2.320740843004842)
    23 = @(2.055961819632387,0.5606415869202124,0.2722823167964208,0.2031497362580978,1.089696314427499,0.07973214163530429,0.3938587522606642,0.2033034070144737,2.324672490841067)
    24 = @(1.789353000266146,0.4842514169519063,0.2591823975887024,0.1939927077186425,1.087286568044988,0.07973214163530429,0.4085686590790001,0.1945939560334722,2.343528173659593)
    25 = @(1.500833966781329,0.4014944072807225,0.2452983013397585,0.184425557943861,1.088010284904712,0.07973214163530429,0.4261087158098569,0.1856208368360797,2.372277102684166)
}

$cols = @(2,3,4,5,6,8,9,10,15)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}
